$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.249.94'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.997.06'
$ws.Range('E3').Value = '  +6.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7754'
$ws.Range('E5').Value = '  +63.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '254.91'
$ws.Range('E6').Value = '  +3.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3483'
$ws.Range('E8').Value = '  +20.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.93'
$ws.Range('E9').Value = '  +25.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07078'
$ws.Range('E10').Value = '  +8.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8445'
$ws.Range('E11').Value = '  +10.59%  '
$ws.Range('E12').Value = '  +5.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '100.95'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.998.29'
$ws.Range('E14').Value = '  +6.18%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.642'
$ws.Range('E15').Value = '  +7.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.26'
$ws.Range('E16').Value = '  +15.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '272.85'
$ws.Range('E17').Value = '  -3.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '31.252.09'
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.983'
$ws.Range('E19').Value = '  +11.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000008023'
$ws.Range('E20').Value = '  +6.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.259.19'
$ws.Range('E21').Value = '  +6.52%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9994'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.106'
$ws.Range('E24').Value = '  +10.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.02'
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.89'
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1427'
$ws.Range('E27').Value = '  +46.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.93'
$ws.Range('E28').Value = '  +4.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.393'
$ws.Range('E29').Value = '  +25.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.598'
$ws.Range('E30').Value = '  +6.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.368'
$ws.Range('E31').Value = '  +2.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.607'
$ws.Range('E32').Value = '  +8.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.449'
$ws.Range('E33').Value = '  +6.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05303'
$ws.Range('E34').Value = '  +9.34%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.248'
$ws.Range('E35').Value = '  +10.27%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7905'
$ws.Range('E36').Value = '  +13.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.769'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02003'
$ws.Range('E38').Value = '  +5.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.915'
$ws.Range('E39').Value = '  +1.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '81.62'
$ws.Range('E40').Value = '  +8.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.775'
$ws.Range('E41').Value = '  +7.39%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.147'
$ws.Range('E42').Value = '  +8.66%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4691'
$ws.Range('E43').Value = '  +10.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8563'
$ws.Range('E44').Value = '  +2.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '105.42'
$ws.Range('E45').Value = '  +3.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.742'
$ws.Range('E47').Value = '  +10.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.01'
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '37.44'
$ws.Range('E49').Value = '  +6.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.567'
$ws.Range('E50').Value = '  +16.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4326'
$ws.Range('E51').Value = '  +9.32%  '
